# The deck currently uses the "Integral" (Red Violet) design on its single
# slide master / theme (ppt/theme/theme1.xml). The author switched the
# presentation's design back to the plain default "Office Theme" palette
# (Design tab -> Office Theme), which re-colors the live theme that backs
# every slide.
#
# Re-apply that by rewriting the active color scheme's twelve theme colors
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the standard Office
# Theme values. The RGB() long used by PowerPoint COM packs color bytes as
# 0xBBGGRR, so each target "RRGGBB" hex is byte-swapped below.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colors = $master.ColorScheme

$colors.Colors(1).RGB  = 0x000000  # dk1      -> 000000
$colors.Colors(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$colors.Colors(3).RGB  = 0x6A5444  # dk2      -> 44546A
$colors.Colors(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$colors.Colors(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$colors.Colors(6).RGB  = 0x317DED  # accent2  -> ED7D31
$colors.Colors(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$colors.Colors(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$colors.Colors(9).RGB  = 0xC47244  # accent5  -> 4472C4
$colors.Colors(10).RGB = 0x47AD70  # accent6  -> 70AD47
$colors.Colors(11).RGB = 0xC16305  # hlink    -> 0563C1
$colors.Colors(12).RGB = 0x724F95  # folHlink -> 954F72
